# PAS-12877: small fix for nightly
# The PartialMatch_CA_SELECT.xlsx test fixture had two extra, unused columns
# (STAT and CHOICE_TIER) inserted in the header/data row. Remove them so the
# sheet matches the expected upload layout again.
#
# Columns as they exist before this edit:
#   ... Y(ANTITHEFTCODE_TEXT) Z(STAT) AA(COLL_SYMBOL) AB(COMP_SYMBOL)
#       AC(CHOICE_TIER) AD(ALTFUEL) AE(BI_SYMBOL) ...
#
# We delete column AC (CHOICE_TIER) first, then column Z (STAT), so that the
# column letters we reference never shift out from under us mid-script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete CHOICE_TIER (currently column AC), then STAT (currently column Z).
# Deleting right-to-left keeps the remaining column references stable.
$ws.Columns("AC:AC").Delete()
$ws.Columns("Z:Z").Delete()

# Restore (as closely as possible) the user's on-screen selection: the two
# now-deleted columns' former positions, which after the shift are occupied
# by COLL_SYMBOL (Z) and BI_SYMBOL (AC).
$ws.Range("Z1:Z1048576").Select()
